# This edit re-orders the 14 species-observation records currently sitting in
# rows 9-22 of the "Artfynd" sheet (each record's id/species/coordinate data
# moves as a unit to a different row), rounds the Ost/Nord (Q/R) coordinates
# to whole metres, and drops the now-unused Starttid/Sluttid (Z/AB) columns.
#
# Rather than hard-coding the shuffled values, we snapshot every relevant
# cell for rows 9-22 live off the open worksheet, compute the new row for
# each record, then write everything back in one pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 9
$lastRow = 22

# Columns (by 1-based index) that participate in this record-level shuffle.
$colMap = @{
    'A'  = 1
    'B'  = 2
    'D'  = 4
    'E'  = 5
    'F'  = 6
    'G'  = 7
    'H'  = 8
    'Q'  = 17
    'R'  = 18
    'AC' = 29
    'AJ' = 36
    'AK' = 37
    'AO' = 41
}

# Row r (old) -> row (new) that record now belongs to.
$rowDestination = @{
    9  = 12
    10 = 18
    11 = 15
    12 = 20
    13 = 9
    14 = 21
    15 = 14
    16 = 13
    17 = 22
    18 = 10
    19 = 17
    20 = 11
    21 = 19
    22 = 16
}

# --- 1. Snapshot every old row's relevant cell values before overwriting anything ---
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($colName in $colMap.Keys) {
        $colIdx = $colMap[$colName]
        $rowData[$colName] = $ws.Cells.Item($r, $colIdx).Value2
    }
    $snapshot[$r] = $rowData
}

# --- 2. Write each record's data into its destination row ---
foreach ($oldRow in $snapshot.Keys) {
    $newRow = $rowDestination[$oldRow]
    $rowData = $snapshot[$oldRow]

    foreach ($colName in $colMap.Keys) {
        $colIdx = $colMap[$colName]
        $value = $rowData[$colName]

        if ($colName -eq 'Q' -or $colName -eq 'R') {
            if ($null -ne $value) {
                $value = [Math]::Round([double]$value, 0)
            }
        }

        $ws.Cells.Item($newRow, $colIdx).Value = $value
    }
}

# --- 3. Drop the Starttid (Z) / Sluttid (AB) columns for every record row ---
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 26).ClearContents() | Out-Null  # Z
    $ws.Cells.Item($r, 28).ClearContents() | Out-Null  # AB
}
